$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename strategy labels / acronyms (column B, rows 2-8)
$ws.Range("B2").Value = "DB Search"
$ws.Range("B3").Value = "SB Search (BS*FS)"
$ws.Range("B4").Value = "DB Search + BS*FS"
$ws.Range("B5").Value = "Scopus + BS*FS"
$ws.Range("B6").Value = "Scopus + BS||FS"
$ws.Range("B7").Value = "Scopus + BS+FS"
$ws.Range("B8").Value = "Scopus + FS+BS"

# Swap the numeric data between row 3 and row 4 (columns C, D, E, I, J, K,
# L, M, N - F, G, H are 0 in both rows so are left untouched), reflecting
# the reordering of the strategies
$c3 = $ws.Cells.Item(3, 3).Value2
$d3 = $ws.Cells.Item(3, 4).Value2
$e3 = $ws.Cells.Item(3, 5).Value2
$i3 = $ws.Cells.Item(3, 9).Value2
$j3 = $ws.Cells.Item(3, 10).Value2
$k3 = $ws.Cells.Item(3, 11).Value2
$l3 = $ws.Cells.Item(3, 12).Value2
$m3 = $ws.Cells.Item(3, 13).Value2
$n3 = $ws.Cells.Item(3, 14).Value2

$c4 = $ws.Cells.Item(4, 3).Value2
$d4 = $ws.Cells.Item(4, 4).Value2
$e4 = $ws.Cells.Item(4, 5).Value2
$i4 = $ws.Cells.Item(4, 9).Value2
$j4 = $ws.Cells.Item(4, 10).Value2
$k4 = $ws.Cells.Item(4, 11).Value2
$l4 = $ws.Cells.Item(4, 12).Value2
$m4 = $ws.Cells.Item(4, 13).Value2
$n4 = $ws.Cells.Item(4, 14).Value2

$ws.Cells.Item(3, 3).Value2 = $c4
$ws.Cells.Item(3, 4).Value2 = $d4
$ws.Cells.Item(3, 5).Value2 = $e4
$ws.Cells.Item(3, 9).Value2 = $i4
$ws.Cells.Item(3, 10).Value2 = $j4
$ws.Cells.Item(3, 11).Value2 = $k4
$ws.Cells.Item(3, 12).Value2 = $l4
$ws.Cells.Item(3, 13).Value2 = $m4
$ws.Cells.Item(3, 14).Value2 = $n4

$ws.Cells.Item(4, 3).Value2 = $c3
$ws.Cells.Item(4, 4).Value2 = $d3
$ws.Cells.Item(4, 5).Value2 = $e3
$ws.Cells.Item(4, 9).Value2 = $i3
$ws.Cells.Item(4, 10).Value2 = $j3
$ws.Cells.Item(4, 11).Value2 = $k3
$ws.Cells.Item(4, 12).Value2 = $l3
$ws.Cells.Item(4, 13).Value2 = $m3
$ws.Cells.Item(4, 14).Value2 = $n3
